$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2
$ws.Range("D2").Value = "57.828.55"
$ws.Range("E2").Value = "  +2.72%  "

# Row 3
$ws.Range("D3").Value = "3.050.56"
$ws.Range("E3").Value = "  +2.49%  "

# Row 4
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
Set-TextCell $ws "D5" "524.84"
$ws.Range("E5").Value = "  +6.15%  "

# Row 6
Set-TextCell $ws "D6" "142.14"
$ws.Range("E6").Value = "  +5.28%  "

# Row 7
$ws.Range("E7").Value = "  +0.09%  "

# Row 8
Set-TextCell $ws "D8" "0.446"
$ws.Range("E8").Value = "  +5.16%  "

# Row 9
Set-TextCell $ws "D9" "7.62"
$ws.Range("E9").Value = "  +4.72%  "

# Row 10
$ws.Range("E10").Value = "  +8.38%  "

# Row 11
Set-TextCell $ws "D11" "0.370"
$ws.Range("E11").Value = "  +5.31%  "

# Row 12
$ws.Range("E12").Value = "  +2.36%  "

# Row 13
$ws.Range("D13").Value = "3.572.61"
$ws.Range("E13").Value = "  +2.52%  "

# Row 14
Set-TextCell $ws "D14" "26.82"
$ws.Range("E14").Value = "  +7.89%  "

# Row 15
Set-TextCell $ws "D15" "0.0000170"
$ws.Range("E15").Value = "  +16.55%  "

# Row 16
$ws.Range("D16").Value = "57.826.66"
$ws.Range("E16").Value = "  +3.08%  "

# Row 17
Set-TextCell $ws "D17" "6.26"
$ws.Range("E17").Value = "  +7.91%  "

# Row 18
$ws.Range("D18").Value = "3.055.02"
$ws.Range("E18").Value = "  +2.72%  "

# Row 19
Set-TextCell $ws "D19" "13.05"
$ws.Range("E19").Value = "  +5.92%  "

# Row 20
Set-TextCell $ws "D20" "8.19"
$ws.Range("E20").Value = "  +6.17%  "

# Row 21
Set-TextCell $ws "D21" "340.05"
$ws.Range("E21").Value = "  +4.89%  "

# Row 22
Set-TextCell $ws "D22" "0.999"
$ws.Range("E22").Value = "  -0.28%  "

# Row 23
Set-TextCell $ws "D23" "5.69"
$ws.Range("E23").Value = "  -0.82%  "

# Row 24
$ws.Range("E24").Value = "  +7.96%  "

# Row 26
$ws.Range("E26").Value = "  +7.48%  "

# Row 27
$ws.Range("D27").Value = "0.0₃0971"
$ws.Range("E27").Value = "  +8.72%  "

# Row 28
$ws.Range("E28").Value = "  -0.04%  "

# Row 29
Set-TextCell $ws "D29" "6.98"
$ws.Range("E29").Value = "  +8.08%  "

# Row 30
Set-TextCell $ws "D30" "7.37"
$ws.Range("E30").Value = "  +9.71%  "

# Row 31
$ws.Range("E31").Value = "  +7.72%  "

# Row 32
$ws.Range("E32").Value = "  +6.67%  "

# Row 33
Set-TextCell $ws "D33" "21.11"
$ws.Range("E33").Value = "  +5.41%  "

# Row 34
Set-TextCell $ws "D34" "4.76"
$ws.Range("E34").Value = "  +6.48%  "

# Row 35
Set-TextCell $ws "D35" "156.50"
$ws.Range("E35").Value = "  +1.05%  "

# Row 36
Set-TextCell $ws "D36" "5.93"
$ws.Range("E36").Value = "  +6.32%  "

# Row 37
Set-TextCell $ws "D37" "1.33"
$ws.Range("E37").Value = "  +4.13%  "

# Row 38
Set-TextCell $ws "D38" "26.00"
$ws.Range("E38").Value = "  +11.62%  "

# Row 39
Set-TextCell $ws "D39" "0.0706"
$ws.Range("E39").Value = "  +4.47%  "

# Row 40
$ws.Range("D40").Value = "3.086.89"
$ws.Range("E40").Value = "  +2.72%  "

# Row 41
Set-TextCell $ws "D41" "37.72"
$ws.Range("E41").Value = "  +3.55%  "

# Row 42
Set-TextCell $ws "D42" "3.89"
$ws.Range("E42").Value = "  +9.47%  "

# Row 43
$ws.Range("E43").Value = "  +0.34%  "

# Row 44
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell $ws "D44" "1.48"
$ws.Range("E44").Value = "  +5.72%  "

# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.340.55"
$ws.Range("E45").Value = "  +6.26%  "

# Row 46
Set-TextCell $ws "D46" "0.661"
$ws.Range("E46").Value = "  +3.91%  "

# Row 47
$ws.Range("E47").Value = "  +3.56%  "

# Row 48
Set-TextCell $ws "D48" "2.02"
$ws.Range("E48").Value = "  +4.64%  "

# Row 49
$ws.Range("E49").Value = "  +4.49%  "

# Row 50
Set-TextCell $ws "D50" "6.05"
$ws.Range("E50").Value = "  +5.35%  "

# Row 51
Set-TextCell $ws "D51" "20.26"
$ws.Range("E51").Value = "  +6.00%  "
